# migrate remote_access import/export to v3, updated unit test
#
# The "NODE" column (D) becomes "NODE_ID" and its values switch from the
# string "master" to the numeric node id 1 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: NODE -> NODE_ID
$ws.Range("D1").Value = "NODE_ID"

# Data rows: node name "master" -> numeric node id 1
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1

# Update the persisted selection / active cell to F3, as in the edited file.
[void]$ws.Range("F3").Select()
